$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update claim data on row 2 (FechaSiniestro, HoraSiniestro, NroPoliza)
# Leading apostrophe keeps these as text values (matching the original
# quotePrefix-based text cells) instead of letting Excel auto-convert
# them to numbers/dates.
$ws.Range("G2").Value = "'22/06/2021"
$ws.Range("H2").Value = "'12:00"
$ws.Range("E2").Value = "'12112002429"

# Update the active selection to F2
$ws.Range("F2").Select()
